# Add new "where" question-tag rows (318-322) to the training data sheet,
# right after the existing data that ends at row 316 (row 317 is a
# deliberate blank separator row, matching the pattern used elsewhere in
# the sheet between tag groups).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("where", "where are you?"),
    @("where", "Where are you?"),
    @("where", "where are you from?"),
    @("where", "Where are u from"),
    @("where", "where are u")
)

$startRow = 318
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
